# Applies the changes described by the commit diff:
#  1. Remove the leftover "_GoBack" bookmark wrapping the document title.
#  2. Delete the entire "TALLER" (workshop) section that had been appended
#     at the end of the document body, right before the final sectPr.

$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark. ---
# "_GoBack" is Word's "last edit position" bookmark; it is hidden from
# Bookmarks.Count / enumeration by default, but it can still be reached
# and removed by name, just like in real Word.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Delete the "TALLER" workshop section at the end of the doc. ---
# Find the paragraph whose text is exactly "TALLER" (the section heading)
# and remove everything from there through the end of the document body
# (i.e. up to, but not including, the final sectPr).
$tallerIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.Trim([char]13, [char]7, [char]10, [char]12)
    if ($t -eq "TALLER") {
        $tallerIndex = $i
    }
}

if ($tallerIndex -gt 0) {
    $sectionStart = $d.Paragraphs($tallerIndex).Range.Start
    $sectionEnd = $d.Content.End
    $d.Range($sectionStart, $sectionEnd).Delete()
}
